$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9-20: column A = req/enc/opt tier, column B = field name, column C = note
# Rows 9-13 are "enc" (encouraged), rows 14-20 are "opt" (optional)

$rows = @(
    @{ r = 9;  a = "enc"; b = "ccodes";       c = "One or more ISO Alpha-2 two-letter country codes (used purely for locative purposes; does not indicate a relationship)" },
    @{ r = 10; a = "enc"; b = "variants";     c = "One or more name and language variants; can be suffixed with @ + language-script code if available" },
    @{ r = 11; a = "enc"; b = "matches";      c = "One or more URIs for matching record(s) in place name authority resources; includes 'exact' and 'close' matches" },
    @{ r = 12; a = "enc"; b = "types";        c = "One or more terms for place type (contributor's term, e.g. pueblo)" },
    @{ r = 13; a = "enc"; b = "aat_typeid";   c = "One or more AAT integer identifiers from WHG's subset list of 160 place type concepts" },
    @{ r = 14; a = "opt"; b = "parent_name";  c = "A single toponym for a containing place" },
    @{ r = 15; a = "opt"; b = "parent_id";    c = "A single URI to a record for the parent_name above" },
    @{ r = 16; a = "opt"; b = "lon";          c = "Longitude for point geometry; decimal degrees" },
    @{ r = 17; a = "opt"; b = "lat";          c = "Latitude for point geometry; decimal degrees" },
    @{ r = 18; a = "opt"; b = "geowkt";       c = "Any geometry in WKT format; polygons should be simplified for rendering performance, using e.g. a GIS function or https://mapshaper.org/; will supercede lon/lat, if both are supplied" },
    @{ r = 19; a = "opt"; b = "geo_srclabel"; c = "Label or short citation for source of the geometry, e.g. GeoNames" },
    @{ r = 20; a = "opt"; b = "geo_srcid";    c = "URI identifier for the source of the geometry, e.g.  http://www.geonames.org/2950159" }
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.r, 1).Value = $row.a
    $ws.Cells.Item($row.r, 2).Value = $row.b
    $ws.Cells.Item($row.r, 3).Value = $row.c
}

$ws.Range("C20").Select()
